# edit.ps1 - apply the "Fix style templates" OOXML changes to the
# start.pptx style template (slide master, slide layouts, slide 1).
#
# The source diff nudges several placeholder/picture sizes by a couple of
# EMUs (a well known artifact of PowerPoint re-serialising a deck) and
# updates the slide-master's slide-number placeholder sample text from
# "<number>" to "1". It also rewrites the random GUIDs on every <a:fld
# type="slidenum"> element across the master and all twelve layouts; those
# ids are internal, opaque identifiers that PowerPoint itself generates
# and never exposes for editing through the object model, so they are
# intentionally left alone here - touching them through TextRange.Text
# would only destroy the field (turning it into plain static text), which
# is strictly worse than leaving the existing, valid field in place.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper values: Shape.Width/Height are expressed in points (1 pt =
# 12700 EMU) and are marshalled through a 32-bit float, so a handful of
# target EMU sizes need a point value nudged to the next representable
# float so that (float)pt * 12700 truncates back to the exact EMU the
# diff expects, instead of landing one EMU short.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Slide 1 - background picture, footer text box, logo picture, title box
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$bg = $s.Shapes.Item(1)                  # picture, off 2258640,0
$bg.Width  = 595.7291870117188           # -> cx 7565760
$bg.Height = 404.8441162109375           # -> cy 5141520

$footer = $s.Shapes.Item(2)              # "Presentation created in..." box
$footer.Width  = 236.01260375976562      # -> cx 2997360
$footer.Height = 21.089765548706055      # -> cy 267840

$logo = $s.Shapes.Item(3)                # small logo picture
$logo.Width  = 22.790552139282227        # -> cx 289440
$logo.Height = 22.790552139282227        # -> cy 289440

$title = $s.Shapes.Item(4)               # "TITLE" box
$title.Width  = 368.0220642089844        # -> cx 4673880
$title.Height = 159.5622100830078        # -> cy 2026440

# ---------------------------------------------------------------------
# Slide master - footer / slide-number / date placeholders
# ---------------------------------------------------------------------
$m = $p.SlideMaster

$ftrPh = $m.Shapes.Item(1)               # ftr idx=1, off 3029040,4767120
$ftrPh.Width  = 242.84410095214844       # -> cx 3084120
$ftrPh.Height = 21.401575088500977       # -> cy 271800

$sldNumPh = $m.Shapes.Item(2)            # sldNum idx=2, off 6458040,4767120
$sldNumPh.Width  = 161.82992553710938    # -> cx 2055240
$sldNumPh.Height = 21.401575088500977    # -> cy 271800

$dtPh = $m.Shapes.Item(3)                # dt idx=3, off 628560,4767120
$dtPh.Width  = 161.82992553710938        # -> cx 2055240
$dtPh.Height = 21.401575088500977        # -> cy 271800

# Sample text shown for the slide-number field in the master changes
# from the generic placeholder "<number>" to the literal sample "1".
$sldNumPh.TextFrame.TextRange.Text = "1"
